# edit.ps1 -- reproduce the commit's content edits via PowerPoint COM-interop.
#
# Summary of the target change:
#   1. The cached "datetime1" field text on the Slide Master and on every one
#      of the 11 Slide Layouts changes from 10/11/2023 -> 11/29/2023.
#   2. Slide 2's "Content Placeholder 8" text is split from one run
#      ("Text paragraph in power point") into two runs
#      ("Text paragraph in " + "power point.") -- i.e. the sentence now ends
#      with a period.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the Date placeholder shown on the Slide Master + every layout.
# ---------------------------------------------------------------------------
$m = $ppt.ActivePresentation.SlideMaster

# Slide Master's own "Date Placeholder 3" shape.
$m.Shapes.Item(3).TextFrame.TextRange.Text = "11/29/2023"

# Each Custom Layout has its own Date placeholder shape, at varying shape
# index (depends on how many other placeholders/shapes the layout has).
$dateShapeIndexByLayout = @{
    1  = 4   # Title Slide              -> Date Placeholder 3
    2  = 4   # Title and Content        -> Date Placeholder 3
    3  = 4   # Section Header           -> Date Placeholder 3
    4  = 5   # Two Content              -> Date Placeholder 4
    5  = 7   # Comparison                -> Date Placeholder 6
    6  = 3   # Title Only               -> Date Placeholder 2
    7  = 2   # Blank                    -> Date Placeholder 1
    8  = 5   # Content with Caption     -> Date Placeholder 4
    9  = 5   # Picture with Caption     -> Date Placeholder 4
    10 = 4   # Title and Vertical Text  -> Date Placeholder 3
    11 = 4   # Vertical Title and Text  -> Date Placeholder 3
}

for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $layout = $m.CustomLayouts.Item($li)
    $dateIdx = $dateShapeIndexByLayout[$li]
    if ($dateIdx) {
        $layout.Shapes.Item($dateIdx).TextFrame.TextRange.Text = "11/29/2023"
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 2: split "Text paragraph in power point" into two runs, adding a
#    trailing period: "Text paragraph in " + "power point."
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$contentShape = $slide2.Shapes.Item(6)   # "Content Placeholder 8"
$textRange = $contentShape.TextFrame.TextRange

$textRange.Text = "Text paragraph in power point."

# Re-touch just the second half of the sentence so it is written back out as
# its own run (separate from the first run), matching the two-run structure
# of the edited deck.
$secondRun = $textRange.Characters(19, 12)
$secondRun.Text = "power point."
